$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 98, shifting existing rows 98:193 down to 99:194
$ws.Rows("98").Insert()

# Populate the newly inserted row 98 with the new record's data
$ws.Range("A98").Value = 7
$ws.Range("B98").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C98").Value = "Ñuble"
$ws.Range("D98").Value = 44539
$ws.Range("E98").Value = 16
$ws.Range("F98").Value = 100112002
$ws.Range("G98").Value = "Pimiento"
$ws.Range("H98").Value = "Cuatro cascos verde"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 160
$ws.Range("K98").Value = 15000
$ws.Range("L98").Value = 15500
$ws.Range("M98").Value = 15250
$ws.Range("N98").Value = "$/caja 15 kilos"
$ws.Range("O98").Value = "Región del Maule"
$ws.Range("P98").Value = 1017
$ws.Range("Q98").Value = 15
$ws.Range("R98").Value = "Hortaliza"
